# regen sval data to filter save games
# Updates computed B-G numeric columns (rows 2-14) on the active sheet
# to the regenerated values (G = B + C + D + E).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3.230985683306322
$ws.Range("C2").Value = 1.667794583268128
$ws.Range("D2").Value = 0.8054896365839992
$ws.Range("E2").Value = 0.496779210170732
$ws.Range("G2").Value = 6.201049113329182
$ws.Range("B3").Value = 3.230985683306322
$ws.Range("C3").Value = 1.667794583268128
$ws.Range("D3").Value = 3.900430680208489
$ws.Range("E3").Value = 0.496779210170732
$ws.Range("G3").Value = 9.295990156953671
$ws.Range("B4").Value = 1.459612070389937
$ws.Range("C4").Value = 1.667794583268128
$ws.Range("D4").Value = 0.1575252929769615
$ws.Range("E4").Value = 0.496779210170732
$ws.Range("G4").Value = 3.781711156805759
$ws.Range("B5").Value = 3.230985683306322
$ws.Range("C5").Value = 1.667794583268128
$ws.Range("D5").Value = 3.900430680208489
$ws.Range("E5").Value = 0.496779210170732
$ws.Range("G5").Value = 9.295990156953671
$ws.Range("B6").Value = 1.459612070389937
$ws.Range("C6").Value = 1.667794583268128
$ws.Range("D6").Value = 0.1575252929769615
$ws.Range("E6").Value = 0.496779210170732
$ws.Range("G6").Value = 3.781711156805759
$ws.Range("B7").Value = 0.6753301551942219
$ws.Range("C7").Value = 1.667794583268128
$ws.Range("D7").Value = 3.900430680208489
$ws.Range("E7").Value = 0.496779210170732
$ws.Range("G7").Value = 6.740334628841572
$ws.Range("B8").Value = 1.459612070389937
$ws.Range("C8").Value = 1.667794583268128
$ws.Range("D8").Value = 3.900430680208489
$ws.Range("E8").Value = 0.496779210170732
$ws.Range("G8").Value = 7.524616544037286
$ws.Range("B9").Value = 0.6753301551942219
$ws.Range("C9").Value = 0.3127903958511391
$ws.Range("D9").Value = 0.8054896365839992
$ws.Range("E9").Value = 0.496779210170732
$ws.Range("G9").Value = 2.290389397800092
$ws.Range("B10").Value = 1.459612070389937
$ws.Range("C10").Value = 1.667794583268128
$ws.Range("D10").Value = 0.8054896365839992
$ws.Range("E10").Value = 0.496779210170732
$ws.Range("G10").Value = 4.429675500412797
$ws.Range("B11").Value = 0.6753301551942219
$ws.Range("C11").Value = 0.3127903958511391
$ws.Range("D11").Value = 0.1575252929769615
$ws.Range("E11").Value = 8.660232485948974
$ws.Range("G11").Value = 9.805878329971296
$ws.Range("B12").Value = 0.6753301551942219
$ws.Range("C12").Value = 0.04240448674262143
$ws.Range("D12").Value = 0.8054896365839992
$ws.Range("E12").Value = 8.660232485948974
$ws.Range("G12").Value = 10.18345676446982
$ws.Range("B13").Value = 3.230985683306322
$ws.Range("C13").Value = 1.667794583268128
$ws.Range("D13").Value = 0.8054896365839992
$ws.Range("E13").Value = 0.496779210170732
$ws.Range("G13").Value = 6.201049113329182
$ws.Range("B14").Value = 3.230985683306322
$ws.Range("C14").Value = 1.667794583268128
$ws.Range("D14").Value = 0.1575252929769615
$ws.Range("E14").Value = 0.496779210170732
$ws.Range("G14").Value = 5.553084769722144
Write-Output "updated sval rows 2-14"
